$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Structural changes -----------------------------------------------
# The trailing column-O cell (always blank) goes away entirely once the
# sheet's used range shrinks back to A:N - clear it before shuffling rows.
$ws.Range("O4").ClearContents()

# Drop the old row 3 ("黃xx" patient record) - it is removed from the log.
$ws.Range("A3:N3").EntireRow.Delete()

# Insert two new rows above the old row 2 ("陳xx" patient) to make room for
# two new patient records ("盧xx" and "吳xx") ahead of the existing data.
$ws.Range("A2:N3").EntireRow.Insert()

# Helper ranges: every column in this sheet stores plain text (even the
# numeric-looking IDs/ages/day-counts), so each of these is formatted as
# text before typing, then reset so no stray formatting lingers on the
# cell once the literal text value is in place.
# (NumberFormat must be applied to one contiguous range at a time - a
# multi-area/union range only formats its first sub-area.)

# --- New row 2 data: 盧xx ----------------------------------------------
$ws.Range("B2:C2").NumberFormat = "@"
$ws.Range("F2").NumberFormat = "@"
$ws.Range("H2").NumberFormat = "@"
$ws.Range("K2").NumberFormat = "@"

$ws.Range("B2").Value = "32"
$ws.Range("C2").Value = "1"
$ws.Range("D2").Value = "盧xx"
$ws.Range("E2").Value = "09C_16_01"
$ws.Range("F2").Value = "3234567"
$ws.Range("G2").Value = "M"
$ws.Range("H2").Value = "79"
$ws.Range("I2").Value = "Colon lesion"
$ws.Range("J2").Value = " right hemicolectomy (lap..."
$ws.Range("K2").Value = "3"
$ws.Range("L2").Value = "王曼玲"
$ws.Range("M2").Value = "15/4h"

$ws.Range("B2:C2").ClearFormats()
$ws.Range("F2").ClearFormats()
$ws.Range("H2").ClearFormats()
$ws.Range("K2").ClearFormats()

# --- New row 3 data: 吳xx ----------------------------------------------
$ws.Range("B3:C3").NumberFormat = "@"
$ws.Range("F3").NumberFormat = "@"
$ws.Range("H3").NumberFormat = "@"
$ws.Range("K3").NumberFormat = "@"

$ws.Range("A3").Value = "PCEA"
$ws.Range("B3").Value = "36"
$ws.Range("C3").Value = "1"
$ws.Range("D3").Value = "吳xx"
$ws.Range("E3").Value = "09B_09_02"
$ws.Range("F3").Value = "4234567"
$ws.Range("G3").Value = "M"
$ws.Range("H3").Value = "72"
$ws.Range("I3").Value = "Liver tumor"
$ws.Range("J3").Value = "Extended right lobectomy"
$ws.Range("K3").Value = "3"
$ws.Range("L3").Value = "王曼玲"
$ws.Range("M3").Value = "4|3|20'|35"
$ws.Range("N3").Value = "mar:265mg+fen:0.5mg /400ML  "

$ws.Range("B3:C3").ClearFormats()
$ws.Range("F3").ClearFormats()
$ws.Range("H3").ClearFormats()
$ws.Range("K3").ClearFormats()

# --- Existing row (old row 4, now row 5): update patient number / age --
$ws.Range("F5").NumberFormat = "@"
$ws.Range("H5").NumberFormat = "@"

$ws.Range("F5").Value = "789"
$ws.Range("H5").Value = "456"

$ws.Range("F5").ClearFormats()
$ws.Range("H5").ClearFormats()
